$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column P ("N de Licencia") and shift existing
# columns (N de Licencia, etc.) one position to the right, matching the
# author's "avances para la carga de trazas" edit that adds an
# "Organismo certificador" field to the DATOS sheet.
$ws.Range("P1").EntireColumn.Insert()

# Populate the header of the newly inserted column.
$ws.Range("P1").Value = "Organismo certificador"
